# "Web 120 / Quiz 03"
# Adds a new "Q03" quiz column (H) to the scores sheet: a header label in H1
# and per-student scores in H2:H16. Column H already existed (formatted,
# just empty) so only the header text + raw score values need to be written;
# all the downstream totals/averages in B17 and B22:F36 are formulas and
# recalculate automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column header (becomes a new shared string "Q03")
$ws.Range("H1").Value = "Q03"

# Per-student Q03 scores
$ws.Range("H2").Value  = 4.5
$ws.Range("H3").Value  = 8
$ws.Range("H4").Value  = 8.5
$ws.Range("H5").Value  = 5
$ws.Range("H6").Value  = 7.5
$ws.Range("H7").Value  = 8
$ws.Range("H8").Value  = 0
$ws.Range("H9").Value  = 8.5
$ws.Range("H10").Value = 7
$ws.Range("H11").Value = 5
$ws.Range("H12").Value = 5.5
$ws.Range("H13").Value = 8
$ws.Range("H14").Value = 4.5
$ws.Range("H15").Value = 5.5
$ws.Range("H16").Value = 6.5

# H11:H16 previously had the default/general alignment (unlike H7:H10 which
# were already centered); bring them in line with the rest of the column's
# number formatting now that they hold real data (also covers the blank
# H17 total-row cell, which picks up the same centered look).
$ws.Range("H11:H17").HorizontalAlignment = -4108

# Selection ends up on B37 after entering the data.
$ws.Range("B37").Select()
